# Generate Report for Handoff
# Updates the status/date/error-detail for the 6cf9b307-...md file across
# the Overview, zh-cn and de-de worksheets, reflecting that the handback
# file is out of date and the item is ready for handoff again.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$overviewDate = "2016-08-18 16:46:57"
$zhDate = "2016-08-18 16:46:51"
$deDate = "2016-08-18 16:46:57"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c59b7159f7fc15bab26546c4209825280fe89a56/e2e/6cf9b307-a186-4fc1-8ccf-4dd66241e2c5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40f54c7261e6da8e857be65127c3ba1e345dbe24/e2e/6cf9b307-a186-4fc1-8ccf-4dd66241e2c5.md."

# --- Overview sheet: row 3 is the 6cf9b307-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $overviewDate

# --- zh-cn sheet: row 3 is the 6cf9b307-...md file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $readyStatus
$wsZh.Range("H3").Value = $zhDate
$wsZh.Range("P3").Value = $errorDetail
# 39.15 is the ColumnWidth (character units) that serializes to the OOXML
# <col width="40"/> the target workbook expects (Excel stores width in its
# own padded/quantized units, so asking for 39.15 lands on 40 after the
# round-trip, the same way the original authoring tool produced it).
$wsZh.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: row 3 is the 6cf9b307-...md file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $readyStatus
$wsDe.Range("H3").Value = $deDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.15
